$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Update sheet1 ("loginTestData"): row 6 password becomes "Abdo@1357" with
#    a mailto hyperlink + the Hyperlink cell style, and row 7 gets a new
#    expectedResult cell.
# ---------------------------------------------------------------------------
$wsLogin = $wb.Worksheets.Item("loginTestData")

$wsLogin.Range("B6").Value = "Abdo@1357"
$null = $wsLogin.Hyperlinks.Add($wsLogin.Range("B6"), "mailto:Abdo@1357")
$wsLogin.Range("B6").Style = "Hyperlink"

$wsLogin.Range("C7").Value = "fail"

# ---------------------------------------------------------------------------
# 2. Add the "googleLogin" sheet right after loginTestData.
# ---------------------------------------------------------------------------
$afterLogin = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsGoogleLogin = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $afterLogin)
$wsGoogleLogin.Name = "googleLogin"

$wsGoogleLogin.Range("A1").Value = "email"
$wsGoogleLogin.Range("B1").Value = "password"
$wsGoogleLogin.Range("A2").Value = "am0643794@gmail.com"
$wsGoogleLogin.Range("B2").Value = "abdo@1357"

$null = $wsGoogleLogin.Hyperlinks.Add($wsGoogleLogin.Range("A2"), "mailto:am0643794@gmail.com")
$wsGoogleLogin.Range("A2").Style = "Hyperlink"

$null = $wsGoogleLogin.Hyperlinks.Add($wsGoogleLogin.Range("B2"), "mailto:abdo@1357")
$wsGoogleLogin.Range("B2").Style = "Hyperlink"

$null = $wsGoogleLogin.Range("F6:G6").Select()

# ---------------------------------------------------------------------------
# 3. Add the "googleForgetUserNamePassword" sheet right after googleLogin.
# ---------------------------------------------------------------------------
$afterGoogleLogin = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForget = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $afterGoogleLogin)
$wsForget.Name = "googleForgetUserNamePassword"

$wsForget.Range("A1").Value = "email"
$wsForget.Range("B1").Value = "password"
$wsForget.Range("C1").Value = "username"
$wsForget.Range("A2").Value = "testmohamed113@gmail.com"
$wsForget.Range("B2").Value = "zskwnnrhdeoxlooo"
$wsForget.Range("C2").Value = "castfer"

$null = $wsForget.Hyperlinks.Add($wsForget.Range("A2"), "mailto:testmohamed113@gmail.com")
$wsForget.Range("A2").Style = "Hyperlink"
$wsForget.Range("B2").Style = "Hyperlink"

$wsForget.PageSetup.Orientation = 1

$null = $wsForget.Range("C2").Select()

# ---------------------------------------------------------------------------
# 4. Final selections / active sheet state: loginTestData keeps a leftover
#    B25 selection, and googleForgetUserNamePassword ends up the active tab.
# ---------------------------------------------------------------------------
$null = $wsLogin.Range("B25").Select()

$null = $wsForget.Activate()
